$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employee")

# G2: email address
$ws.Range("G2").Value = "koshtech.site@gmail.com"

# O2, P2, Q2, U2: set numeric value 1 (previously empty inline strings)
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("U2").Value = 1
